$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 361; everything from old row 361 downward
# shifts down by one (old 361 -> 362, old 402 -> 403, old 403 -> 404).
$ws.Rows.Item(361).Insert()

# Populate the newly inserted row 361 with the new weekly record.
$ws.Cells.Item(361, 1).Value = 9
$ws.Cells.Item(361, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(361, 3).Value = "Metropolitana"
$ws.Cells.Item(361, 4).Value = 45212
$ws.Cells.Item(361, 5).Value = 13
$ws.Cells.Item(361, 6).Value = 100112026
$ws.Cells.Item(361, 7).Value = "Haba"
$ws.Cells.Item(361, 8).Value = "Sin especificar"
$ws.Cells.Item(361, 9).Value = "Primera"
$ws.Cells.Item(361, 10).Value = 70
$ws.Cells.Item(361, 11).Value = 9000
$ws.Cells.Item(361, 12).Value = 10000
$ws.Cells.Item(361, 13).Value = 9500
$ws.Cells.Item(361, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(361, 15).Value = "Provincia de Melipilla"
$ws.Cells.Item(361, 16).Value = 380
$ws.Cells.Item(361, 17).Value = 25
$ws.Cells.Item(361, 18).Value = "Hortaliza"
